$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 419.18182
$ws.Range("I33").Value = 391.57144
$ws.Range("M33").Value = -162.57144
$ws.Range("K33").Value = 391.57144
# Row 40
$ws.Range("H40").Value = 1022.7907
$ws.Range("I40").Value = 1196
$ws.Range("M40").Value = -1021
$ws.Range("K40").Value = 1196
$ws.Range("N40").Value = -1350
$ws.Range("L40").Value = 1000
$ws.Range("J40").Value = 1000
# Row 45
$ws.Range("H45").Value = 2000
$ws.Range("I45").Value = 1000
$ws.Range("K45").Value = 3000
$ws.Range("M45").Value = -2808
# Row 64
$ws.Range("H64").Value = 4134.4585
$ws.Range("I64").Value = 3790
$ws.Range("M64").Value = -3542
$ws.Range("K64").Value = 3790
$ws.Range("N64").Value = -4876.5
$ws.Range("L64").Value = 4380.5
$ws.Range("J64").Value = 4380.5
# Row 67
$ws.Range("H67").Value = 4134.4585
$ws.Range("I67").Value = 3790
$ws.Range("M67").Value = -2932
$ws.Range("K67").Value = 3790
$ws.Range("N67").Value = -6096.5
$ws.Range("L67").Value = 4380.5
$ws.Range("J67").Value = 4380.5
# Row 74
$ws.Range("H74").Value = 3774.138
$ws.Range("I74").Value = 3612.1177
$ws.Range("M74").Value = -2676.1177
$ws.Range("K74").Value = 3612.1177
$ws.Range("N74").Value = -5875.6667
$ws.Range("L74").Value = 4003.6667
$ws.Range("J74").Value = 4003.6667
# Row 76
$ws.Range("H76").Value = 3599.25
$ws.Range("I76").Value = 3442.7144
$ws.Range("M76").Value = -3127.7144
$ws.Range("K76").Value = 3442.7144
$ws.Range("N76").Value = -4448.4
$ws.Range("L76").Value = 3818.4
$ws.Range("J76").Value = 3818.4
# Row 77
$ws.Range("H77").Value = 3774.138
$ws.Range("I77").Value = 3612.1177
$ws.Range("M77").Value = -13380.5885
$ws.Range("K77").Value = 18060.5885
$ws.Range("N77").Value = -29378.3335
$ws.Range("L77").Value = 20018.3335
$ws.Range("J77").Value = 4003.6667
# Row 79
$ws.Range("H79").Value = 3599.25
$ws.Range("I79").Value = 3442.7144
$ws.Range("M79").Value = -2350.7144
$ws.Range("K79").Value = 3442.7144
$ws.Range("N79").Value = -6002.4
$ws.Range("L79").Value = 3818.4
$ws.Range("J79").Value = 3818.4
# Row 112
$ws.Range("H112").Value = 34485108
$ws.Range("N112").Value = -9680.110999999999
$ws.Range("L112").Value = 7464.110999999999
$ws.Range("J112").Value = 2488.037
# Row 132
$ws.Range("H132").Value = 892803.25
$ws.Range("I132").Value = 1783.6531
$ws.Range("M132").Value = -2820.9593
$ws.Range("K132").Value = 5350.9593
$ws.Range("N132").Value = -24513449
$ws.Range("L132").Value = 24508389
$ws.Range("J132").Value = 8169463

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1611.2433
$ws.Range("I2").Value = 1518.1852
$ws.Range("M2").Value = -1405.1852
$ws.Range("K2").Value = 1518.1852
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("L12").Value = 0
$ws.Range("J12").Value = 0
# Row 23
$ws.Range("H23").Value = 40000
$ws.Range("I23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("K23").Value = 0
$ws.Range("N23").Value = -40518
$ws.Range("L23").Value = 40000
$ws.Range("J23").Value = 40000
# Row 32
$ws.Range("H32").Value = 19383.096
$ws.Range("I32").Value = 18331.027
$ws.Range("M32").Value = -18044.027
$ws.Range("K32").Value = 18331.027
$ws.Range("N32").Value = -23849.75
$ws.Range("L32").Value = 23275.75
$ws.Range("J32").Value = 23275.75
# Row 37
$ws.Range("H37").Value = 8800
$ws.Range("N37").Value = -38546
$ws.Range("L37").Value = 38000
$ws.Range("J37").Value = 38000
# Row 44
$ws.Range("H44").Value = 20049
$ws.Range("N44").Value = -21025
$ws.Range("L44").Value = 20049
$ws.Range("J44").Value = 20049
# Row 55
$ws.Range("H55").Value = 20053
$ws.Range("N55").Value = -20683
$ws.Range("L55").Value = 20053
$ws.Range("J55").Value = 20053
# Row 61
$ws.Range("H61").Value = 12221137
$ws.Range("I61").Value = 13347963
$ws.Range("M61").Value = -13347751
$ws.Range("K61").Value = 13347963
$ws.Range("N61").Value = -148428
$ws.Range("L61").Value = 148004
$ws.Range("J61").Value = 148004
# Row 80
$ws.Range("H80").Value = 34993.332
$ws.Range("N80").Value = -36986
$ws.Range("L80").Value = 34990
$ws.Range("J80").Value = 34990
# Row 83
$ws.Range("H83").Value = 34993.332
$ws.Range("N83").Value = -114954
$ws.Range("L83").Value = 104970
$ws.Range("J83").Value = 34990
# Row 116
$ws.Range("H116").Value = 1611.2433
$ws.Range("I116").Value = 1518.1852
$ws.Range("M116").Value = 775.8148000000001
$ws.Range("K116").Value = 1518.1852
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("L133").Value = 0
$ws.Range("J133").Value = 0
# Row 134
$ws.Range("H134").Value = 54756.617
$ws.Range("N134").Value = -64896.617
$ws.Range("L134").Value = 54756.617
$ws.Range("J134").Value = 54756.617
# Row 136
$ws.Range("H136").Value = 12221137
$ws.Range("I136").Value = 13347963
$ws.Range("M136").Value = -40041339
$ws.Range("K136").Value = 40043889
$ws.Range("N136").Value = -449112
$ws.Range("L136").Value = 444012
$ws.Range("J136").Value = 148004

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1611.2433
$ws.Range("I3").Value = 1518.1852
$ws.Range("M3").Value = -1404.1852
$ws.Range("K3").Value = 1518.1852
# Row 99
$ws.Range("H99").Value = 1030.0588
$ws.Range("I99").Value = 917.1818
$ws.Range("M99").Value = 580.8182
$ws.Range("K99").Value = 917.1818
$ws.Range("N99").Value = -4233
$ws.Range("L99").Value = 1237
$ws.Range("J99").Value = 1237

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 3678.6667
$ws.Range("I62").Value = 3225.7144
$ws.Range("M62").Value = -2601.7144
$ws.Range("K62").Value = 3225.7144
$ws.Range("N62").Value = -5323
$ws.Range("L62").Value = 4075
$ws.Range("J62").Value = 4075
# Row 65
$ws.Range("H65").Value = 3678.6667
$ws.Range("I65").Value = 3225.7144
$ws.Range("M65").Value = -13008.572
$ws.Range("K65").Value = 16128.572
$ws.Range("N65").Value = -26615
$ws.Range("L65").Value = 20375
$ws.Range("J65").Value = 4075
# Row 135
$ws.Range("H135").Value = 32780
$ws.Range("N135").Value = -42920
$ws.Range("L135").Value = 32780
$ws.Range("J135").Value = 32780
# Row 141
$ws.Range("H141").Value = 51330.4
$ws.Range("N141").Value = -62023
$ws.Range("L141").Value = 51663
$ws.Range("J141").Value = 51663

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("L40").Value = 0
$ws.Range("J40").Value = 0
# Row 80
$ws.Range("H80").Value = 3690.111
$ws.Range("I80").Value = 3150
$ws.Range("M80").Value = -2152
$ws.Range("K80").Value = 3150
$ws.Range("N80").Value = -5956.1667
$ws.Range("L80").Value = 3960.1667
$ws.Range("J80").Value = 3960.1667
# Row 83
$ws.Range("H83").Value = 3690.111
$ws.Range("I83").Value = 3150
$ws.Range("M83").Value = -10758
$ws.Range("K83").Value = 15750
$ws.Range("N83").Value = -29784.8335
$ws.Range("L83").Value = 19800.8335
$ws.Range("J83").Value = 3960.1667

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 831.06665
$ws.Range("I82").Value = 815.25
$ws.Range("M82").Value = -454.25
$ws.Range("K82").Value = 815.25
$ws.Range("N82").Value = -1558.8182
$ws.Range("L82").Value = 836.8182
$ws.Range("J82").Value = 836.8182
# Row 85
$ws.Range("H85").Value = 831.06665
$ws.Range("I85").Value = 815.25
$ws.Range("M85").Value = 432.75
$ws.Range("K85").Value = 815.25
$ws.Range("N85").Value = -3332.8182
$ws.Range("L85").Value = 836.8182
$ws.Range("J85").Value = 836.8182
# Row 100
$ws.Range("H100").Value = 1549.8
$ws.Range("I100").Value = 1323.5385
$ws.Range("M100").Value = -782.5385000000001
$ws.Range("K100").Value = 1323.5385
$ws.Range("N100").Value = -3052
$ws.Range("L100").Value = 1970
$ws.Range("J100").Value = 1970
